$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (small)
$ws.Range("C2").Value = 1.676
$ws.Range("D2").Value = 10.055
$ws.Range("E2").Value = 3.926
$ws.Range("F2").Value = 3.95
$ws.Range("G2").Value = 6.56
$ws.Range("H2").Value = 3.974
$ws.Range("I2").Value = 1.245
$ws.Range("J2").Value = 3.735
$ws.Range("K2").Value = 8.02
$ws.Range("L2").Value = 7.972
$ws.Range("M2").Value = 3.495
$ws.Range("N2").Value = 3.184
$ws.Range("O2").Value = 57.792

# Row 3 (medium)
$ws.Range("D3").Value = 2.825
$ws.Range("E3").Value = 1.484
$ws.Range("F3").Value = 1.628
$ws.Range("G3").Value = 2.179
$ws.Range("H3").Value = 1.652
$ws.Range("I3").Value = 0.9340000000000001
$ws.Range("J3").Value = 2.418
$ws.Range("K3").Value = 3.376
$ws.Range("L3").Value = 3.256
$ws.Range("M3").Value = 2.011
$ws.Range("N3").Value = 2.035
$ws.Range("O3").Value = 24.277

# Row 4 (large)
$ws.Range("C4").Value = 0.168
$ws.Range("D4").Value = 1.341
$ws.Range("E4").Value = 0.67
$ws.Range("F4").Value = 0.766
$ws.Range("G4").Value = 0.766
$ws.Range("H4").Value = 1.197
$ws.Range("I4").Value = 0.766
$ws.Range("J4").Value = 1.293
$ws.Range("K4").Value = 2.107
$ws.Range("L4").Value = 1.532
$ws.Range("M4").Value = 1.077
$ws.Range("N4").Value = 1.269
$ws.Range("O4").Value = 12.952

# Row 5 (huge)
$ws.Range("D5").Value = 0.048
$ws.Range("H5").Value = 0.239
$ws.Range("O5").Value = 0.287

# Row 6 (unknown_sz)
$ws.Range("C6").Value = 0.239
$ws.Range("D6").Value = 0.622
$ws.Range("E6").Value = 0.12
$ws.Range("F6").Value = 0.263
$ws.Range("G6").Value = 0.407
$ws.Range("H6").Value = 0.527
$ws.Range("I6").Value = 0.311
$ws.Range("J6").Value = 0.335
$ws.Range("K6").Value = 0.599
$ws.Range("L6").Value = 0.694
$ws.Range("M6").Value = 0.192
$ws.Range("N6").Value = 0.383
$ws.Range("O6").Value = 4.692

# Row 7 (COL_TOT)
$ws.Range("C7").Value = 2.562
$ws.Range("D7").Value = 14.891
$ws.Range("E7").Value = 6.2
$ws.Range("F7").Value = 6.607
$ws.Range("G7").Value = 9.911999999999999
$ws.Range("H7").Value = 7.589
$ws.Range("I7").Value = 3.256
$ws.Range("J7").Value = 7.781
$ws.Range("K7").Value = 14.102
$ws.Range("L7").Value = 13.454
$ws.Range("M7").Value = 6.775
$ws.Range("N7").Value = 6.871
$ws.Range("O7").Value = 100
